$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 553.3333
$ws.Range("J2").Value = 1375
$ws.Range("L2").Value = 1375
$ws.Range("N2").Value = -1601
$ws.Range("H6").Value = 1091.8
$ws.Range("I6").Value = 1197
$ws.Range("J6").Value = 145
$ws.Range("K6").Value = 3591
$ws.Range("L6").Value = 435
$ws.Range("N6").Value = -659
$ws.Range("H9").Value = 8078.846
$ws.Range("I9").Value = 11323.111
$ws.Range("J9").Value = 779.25
$ws.Range("K9").Value = 11323.111
$ws.Range("L9").Value = 779.25
$ws.Range("M9").Value = -11154.111
$ws.Range("N9").Value = -1117.25
$ws.Range("H17").Value = 203950.64
$ws.Range("J17").Value = 213804.58
$ws.Range("L17").Value = 641413.74
$ws.Range("N17").Value = -641749.74
$ws.Range("H18").Value = 1500
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = $null
$ws.Range("H19").Value = 1559
$ws.Range("I19").Value = 1378
$ws.Range("K19").Value = 1378
$ws.Range("M19").Value = -1203
$ws.Range("H29").Value = 5516.5
$ws.Range("J29").Value = 5516.5
$ws.Range("L29").Value = 16549.5
$ws.Range("N29").Value = -17111.5
$ws.Range("H33").Value = 79663.21000000001
$ws.Range("I33").Value = 139029.5
$ws.Range("K33").Value = 139029.5
$ws.Range("M33").Value = -138800.5
$ws.Range("H38").Value = 1739.8462
$ws.Range("I38").Value = 1239.2222
$ws.Range("J38").Value = 2866.25
$ws.Range("K38").Value = 3717.6666
$ws.Range("L38").Value = 8598.75
$ws.Range("M38").Value = -3345.6666
$ws.Range("N38").Value = -9342.75
$ws.Range("H40").Value = 6913.5713
$ws.Range("J40").Value = 8223.75
$ws.Range("L40").Value = 8223.75
$ws.Range("N40").Value = -8573.75
$ws.Range("H43").Value = 9999
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").Value = $null
$ws.Range("H58").Value = 4821.3335
$ws.Range("I58").Value = 976.3333
$ws.Range("J58").Value = 8666.333000000001
$ws.Range("K58").Value = 2928.9999
$ws.Range("L58").Value = 25998.999
$ws.Range("M58").Value = -2778.9999
$ws.Range("N58").Value = -26298.999
$ws.Range("H63").Value = 85899.5
$ws.Range("I63").Value = 85899.5
$ws.Range("K63").Value = 85899.5
$ws.Range("M63").Value = -85275.5
$ws.Range("H66").Value = 85899.5
$ws.Range("I66").Value = 85899.5
$ws.Range("K66").Value = 257698.5
$ws.Range("M66").Value = -254578.5
$ws.Range("H74").Value = 2256747.5
$ws.Range("I74").Value = 2256747.5
$ws.Range("K74").Value = 2256747.5
$ws.Range("M74").Value = -2255811.5
$ws.Range("H75").Value = 29407
$ws.Range("J75").Value = 29407
$ws.Range("L75").Value = 29407
$ws.Range("N75").Value = -31279
$ws.Range("H76").Value = 3182.6667
$ws.Range("J76").Value = 3940
$ws.Range("L76").Value = 3940
$ws.Range("N76").Value = -4570
$ws.Range("H77").Value = 2256747.5
$ws.Range("I77").Value = 2256747.5
$ws.Range("K77").Value = 11283737.5
$ws.Range("M77").Value = -11279057.5
$ws.Range("H78").Value = 29407
$ws.Range("J78").Value = 29407
$ws.Range("L78").Value = 88221
$ws.Range("N78").Value = -97581
$ws.Range("H79").Value = 3182.6667
$ws.Range("J79").Value = 3940
$ws.Range("L79").Value = 3940
$ws.Range("N79").Value = -6124
$ws.Range("H98").Value = 1589.8334
$ws.Range("I98").Value = 1095.1177
$ws.Range("K98").Value = 1095.1177
$ws.Range("M98").Value = 402.8823
$ws.Range("H113").Value = 7215
$ws.Range("I113").Value = 7195
$ws.Range("J113").Value = 7225
$ws.Range("K113").Value = 7195
$ws.Range("L113").Value = 7225
$ws.Range("M113").Value = -3941
$ws.Range("N113").Value = -13733
$ws.Range("H122").Value = 1589.8334
$ws.Range("I122").Value = 1095.1177
$ws.Range("K122").Value = 3285.3531
$ws.Range("M122").Value = -835.3531000000003
$ws.Range("H137").Value = 2754.7273
$ws.Range("J137").Value = 4597
$ws.Range("L137").Value = 13791
$ws.Range("N137").Value = -18891

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 2323.2856
$ws.Range("I14").Value = 210.66667
$ws.Range("J14").Value = 14999
$ws.Range("K14").Value = 210.66667
$ws.Range("L14").Value = 14999
$ws.Range("N14").Value = -15349
$ws.Range("H24").Value = 24250
$ws.Range("J24").Value = 24250
$ws.Range("L24").Value = 24250
$ws.Range("N24").Value = -24998
$ws.Range("H28").Value = 8333.333000000001
$ws.Range("J28").Value = 16001
$ws.Range("L28").Value = 16001
$ws.Range("N28").Value = -16385
$ws.Range("H45").Value = 6827.5713
$ws.Range("I45").Value = 7480.1816
$ws.Range("K45").Value = 7480.1816
$ws.Range("M45").Value = -7103.1816
$ws.Range("H61").Value = 45457530
$ws.Range("I61").Value = 52634596
$ws.Range("K61").Value = 52634596
$ws.Range("M61").Value = -52634384
$ws.Range("H63").Value = 1729.8572
$ws.Range("I63").Value = 1754.5385
$ws.Range("K63").Value = 1754.5385
$ws.Range("M63").Value = -1068.5385
$ws.Range("H64").Value = 71900
$ws.Range("I64").Value = 71900
$ws.Range("K64").Value = 71900
$ws.Range("M64").Value = -71652
$ws.Range("H66").Value = 1729.8572
$ws.Range("I66").Value = 1754.5385
$ws.Range("K66").Value = 8772.692500000001
$ws.Range("M66").Value = -5340.692500000001
$ws.Range("H67").Value = 71900
$ws.Range("I67").Value = 71900
$ws.Range("K67").Value = 71900
$ws.Range("M67").Value = -71042
$ws.Range("H74").Value = 55559704
$ws.Range("I74").Value = 76925890
$ws.Range("J74").Value = 7622.2
$ws.Range("K74").Value = 76925890
$ws.Range("L74").Value = 7622.2
$ws.Range("M74").Value = -76925016
$ws.Range("N74").Value = -9370.200000000001
$ws.Range("H77").Value = 55559704
$ws.Range("I77").Value = 76925890
$ws.Range("J77").Value = 7622.2
$ws.Range("K77").Value = 384629450
$ws.Range("L77").Value = 38111
$ws.Range("M77").Value = -384625082
$ws.Range("N77").Value = -46847
$ws.Range("H94").Value = 19330
$ws.Range("J94").Value = 19330
$ws.Range("N94").Value = -21132
$ws.Range("H98").Value = 26227
$ws.Range("J98").Value = 26227
$ws.Range("L98").Value = 26227
$ws.Range("N98").Value = -32217
$ws.Range("H99").Value = 8333.333000000001
$ws.Range("J99").Value = 16001
$ws.Range("L99").Value = 16001
$ws.Range("N99").Value = -21991
$ws.Range("H100").Value = 24250
$ws.Range("J100").Value = 24250
$ws.Range("L100").Value = 24250
$ws.Range("N100").Value = -26414
$ws.Range("H101").Value = 105414
$ws.Range("J101").Value = 105414
$ws.Range("L101").Value = 105414
$ws.Range("N101").Value = -111904
$ws.Range("H102").Value = 11113133
$ws.Range("I102").Value = 14287456
$ws.Range("K102").Value = 14287456
$ws.Range("M102").Value = -14285834
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("N106").Value = $null
$ws.Range("H110").Value = 112226.445
$ws.Range("I110").Value = 112226.445
$ws.Range("K110").Value = 112226.445
$ws.Range("M110").Value = -110181.445
$ws.Range("H132").Value = 4169657
$ws.Range("I132").Value = 5002783.5
$ws.Range("K132").Value = 15008350.5
$ws.Range("M132").Value = -15005820.5
$ws.Range("H136").Value = 45457530
$ws.Range("I136").Value = 52634596
$ws.Range("K136").Value = 157903788
$ws.Range("M136").Value = -157901238

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3900
$ws.Range("J20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5494
$ws.Range("H22").Value = 1838.35
$ws.Range("J22").Value = 1192.4445
$ws.Range("L22").Value = 1192.4445
$ws.Range("N22").Value = -1538.4445
$ws.Range("H28").Value = 29900
$ws.Range("J28").Value = 29900
$ws.Range("L28").Value = 29900
$ws.Range("N28").Value = -30488
$ws.Range("H37").Value = 2304.1738
$ws.Range("I37").Value = 1272.591
$ws.Range("J37").Value = 24999
$ws.Range("K37").Value = 1272.591
$ws.Range("L37").Value = 24999
$ws.Range("N37").Value = -25273
$ws.Range("H46").Value = 10000
$ws.Range("I46").Value = 10000
$ws.Range("K46").Value = 10000
$ws.Range("M46").Value = -9702
$ws.Range("H94").Value = 6789.5
$ws.Range("I94").Value = 10379.5
$ws.Range("J94").Value = 3199.5
$ws.Range("K94").Value = 10379.5
$ws.Range("L94").Value = 3199.5
$ws.Range("M94").Value = -9928.5
$ws.Range("N94").Value = -4101.5
$ws.Range("H107").Value = 88507.664
$ws.Range("I107").Value = 1121.4445
$ws.Range("J107").Value = 350666.34
$ws.Range("K107").Value = 1121.4445
$ws.Range("L107").Value = 350666.34
$ws.Range("M107").Value = 798.5554999999999
$ws.Range("N107").Value = -354506.34
$ws.Range("H134").Value = 18524690
$ws.Range("I134").Value = 20006306
$ws.Range("K134").Value = 60018918
$ws.Range("M134").Value = -60016383

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 8468.786
$ws.Range("I7").Value = 13115.777
$ws.Range("J7").Value = 104.2
$ws.Range("K7").Value = 13115.777
$ws.Range("L7").Value = 104.2
$ws.Range("M7").Value = -13002.777
$ws.Range("N7").Value = -330.2
$ws.Range("H31").Value = 4931.72
$ws.Range("I31").Value = 3405.4443
$ws.Range("J31").Value = 8856.429
$ws.Range("K31").Value = 3405.4443
$ws.Range("L31").Value = 8856.429
$ws.Range("M31").Value = -3110.4443
$ws.Range("N31").Value = -9446.429
$ws.Range("H34").Value = 4931.72
$ws.Range("I34").Value = 3405.4443
$ws.Range("J34").Value = 8856.429
$ws.Range("K34").Value = 3405.4443
$ws.Range("L34").Value = 8856.429
$ws.Range("M34").Value = -3203.4443
$ws.Range("N34").Value = -9260.429
$ws.Range("H58").Value = 26324852
$ws.Range("I58").Value = 38473450
$ws.Range("K58").Value = 38473450
$ws.Range("M58").Value = -38473247
$ws.Range("H86").Value = 11431
$ws.Range("I86").Value = 8517.25
$ws.Range("J86").Value = 12726
$ws.Range("K86").Value = 8517.25
$ws.Range("L86").Value = 12726
$ws.Range("M86").Value = -7394.25
$ws.Range("N86").Value = -14972
$ws.Range("H89").Value = 11431
$ws.Range("I89").Value = 8517.25
$ws.Range("J89").Value = 12726
$ws.Range("K89").Value = 42586.25
$ws.Range("L89").Value = 63630
$ws.Range("M89").Value = -36970.25
$ws.Range("N89").Value = -74862
$ws.Range("H95").Value = 18295.334
$ws.Range("J95").Value = 18295.334
$ws.Range("L95").Value = 18295.334
$ws.Range("N95").Value = -23787.334
$ws.Range("H96").Value = 13949.167
$ws.Range("J96").Value = 13949.167
$ws.Range("L96").Value = 13949.167
$ws.Range("N96").Value = -19441.167
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3000
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = $null
$ws.Range("H105").Value = 2268810
$ws.Range("I105").Value = 2268810
$ws.Range("K105").Value = 2268810
$ws.Range("M105").Value = -2267063
$ws.Range("H122").Value = 1882.2667
$ws.Range("I122").Value = 2198.9
$ws.Range("J122").Value = 1249
$ws.Range("K122").Value = 6596.700000000001
$ws.Range("L122").Value = 3747
$ws.Range("N122").Value = -8647
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = $null
$ws.Range("H132").Value = 37040204
$ws.Range("I132").Value = 47621976
$ws.Range("K132").Value = 142865928
$ws.Range("M132").Value = -142863398
$ws.Range("H134").Value = 19233954
$ws.Range("I134").Value = 27780334
$ws.Range("K134").Value = 83341002
$ws.Range("M134").Value = -83338467
$ws.Range("H136").Value = 26324852
$ws.Range("I136").Value = 38473450
$ws.Range("K136").Value = 115420350
$ws.Range("M136").Value = -115417800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 189687.5
$ws.Range("J11").Value = 100000
$ws.Range("L11").Value = 300000
$ws.Range("N11").Value = -300280
$ws.Range("H68").Value = 2910.6428
$ws.Range("I68").Value = 2953
$ws.Range("K68").Value = 8859
$ws.Range("M68").Value = -8048
$ws.Range("H70").Value = 14858.4
$ws.Range("I70").Value = 6716.8
$ws.Range("K70").Value = 20150.4
$ws.Range("M70").Value = -19835.4
$ws.Range("H71").Value = 2910.6428
$ws.Range("I71").Value = 2953
$ws.Range("K71").Value = 26577
$ws.Range("M71").Value = -22521
$ws.Range("H73").Value = 14858.4
$ws.Range("I73").Value = 6716.8
$ws.Range("K73").Value = 20150.4
$ws.Range("M73").Value = -19058.4
$ws.Range("H75").Value = 1600
$ws.Range("J75").Value = 1600
$ws.Range("N75").Value = -6796
$ws.Range("H78").Value = 1600
$ws.Range("J78").Value = 1600
$ws.Range("N78").Value = -24384
$ws.Range("H92").Value = 422
$ws.Range("I92").Value = 500
$ws.Range("K92").Value = 1500
$ws.Range("M92").Value = -252
$ws.Range("H103").Value = 327
$ws.Range("I103").Value = 327
$ws.Range("K103").Value = 981
$ws.Range("M103").Value = -102
$ws.Range("H113").Value = 125504.5
$ws.Range("I113").Value = 999999
$ws.Range("J113").Value = 576.7143
$ws.Range("K113").Value = 2999997
$ws.Range("L113").Value = 1730.1429
$ws.Range("M113").Value = -2997827
$ws.Range("N113").Value = -6070.1429
$ws.Range("H131").Value = 2076.3333
$ws.Range("I131").Value = 2076.3333
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 6228.999899999999
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 121
$ws.Range("I2").Value = 65.333336
$ws.Range("K2").Value = 65.333336
$ws.Range("M2").Value = 47.666664
$ws.Range("H3").Value = 396.6
$ws.Range("I3").Value = 470.75
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 470.75
$ws.Range("L3").Value = 100
$ws.Range("M3").Value = -354.75
$ws.Range("N3").Value = -332
$ws.Range("H39").Value = 8001
$ws.Range("J39").Value = 8001
$ws.Range("L39").Value = 8001
$ws.Range("N39").Value = -9065
$ws.Range("H80").Value = 4911.6665
$ws.Range("J80").Value = 4737
$ws.Range("L80").Value = 4737
$ws.Range("N80").Value = -6733
$ws.Range("H83").Value = 4911.6665
$ws.Range("J83").Value = 4737
$ws.Range("L83").Value = 23685
$ws.Range("N83").Value = -33669
$ws.Range("H96").Value = 32990.5
$ws.Range("J96").Value = 32990.5
$ws.Range("L96").Value = 32990.5
$ws.Range("N96").Value = -38482.5
$ws.Range("H99").Value = 9800
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = $null
$ws.Range("H101").Value = 19999
$ws.Range("J101").Value = 19999
$ws.Range("L101").Value = 19999
$ws.Range("N101").Value = -26489
$ws.Range("H107").Value = 3819.8572
$ws.Range("J107").Value = 10002.5
$ws.Range("L107").Value = 10002.5
$ws.Range("N107").Value = -13842.5
$ws.Range("H113").Value = 88856
$ws.Range("J113").Value = 3500
$ws.Range("L113").Value = 3500
$ws.Range("N113").Value = -7840
$ws.Range("H122").Value = 3900.7368
$ws.Range("I122").Value = 3138
$ws.Range("J122").Value = 5553.3335
$ws.Range("K122").Value = 9414
$ws.Range("L122").Value = 16660.0005
$ws.Range("M122").Value = -6964
$ws.Range("N122").Value = -21560.0005
$ws.Range("H132").Value = 9619805
$ws.Range("I132").Value = 12504548
$ws.Range("K132").Value = 37513644
$ws.Range("M132").Value = -37511114

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2973.5557
$ws.Range("I22").Value = 2823.2856
$ws.Range("K22").Value = 2823.2856
$ws.Range("M22").Value = -2528.2856
$ws.Range("H27").Value = 2973.5557
$ws.Range("I27").Value = 2823.2856
$ws.Range("K27").Value = 2823.2856
$ws.Range("M27").Value = -2716.2856
$ws.Range("H61").Value = 749.5
$ws.Range("I61").Value = 749.5
$ws.Range("K61").Value = 749.5
$ws.Range("M61").Value = -547.5
$ws.Range("H68").Value = 1993.25
$ws.Range("I68").Value = 1432.6666
$ws.Range("J68").Value = 2329.6
$ws.Range("K68").Value = 1432.6666
$ws.Range("L68").Value = 2329.6
$ws.Range("M68").Value = -683.6666
$ws.Range("N68").Value = -3827.6
$ws.Range("H71").Value = 1993.25
$ws.Range("I71").Value = 1432.6666
$ws.Range("J71").Value = 2329.6
$ws.Range("K71").Value = 7163.333000000001
$ws.Range("L71").Value = 11648
$ws.Range("M71").Value = -3419.333000000001
$ws.Range("N71").Value = -19136
$ws.Range("H82").Value = 1645.625
$ws.Range("I82").Value = 1532.6666
$ws.Range("J82").Value = 1984.5
$ws.Range("K82").Value = 1532.6666
$ws.Range("L82").Value = 1984.5
$ws.Range("M82").Value = -1171.6666
$ws.Range("N82").Value = -2706.5
$ws.Range("H85").Value = 1645.625
$ws.Range("I85").Value = 1532.6666
$ws.Range("J85").Value = 1984.5
$ws.Range("K85").Value = 1532.6666
$ws.Range("L85").Value = 1984.5
$ws.Range("M85").Value = -284.6666
$ws.Range("N85").Value = -4480.5
$ws.Range("H113").Value = 749.5
$ws.Range("I113").Value = 749.5
$ws.Range("K113").Value = 749.5
$ws.Range("M113").Value = 1420.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 100000
$ws.Range("J7").Value = 100000
$ws.Range("N7").Value = -100226
$ws.Range("H14").Value = 1234
$ws.Range("I14").Value = 1234
$ws.Range("K14").Value = 1234
$ws.Range("M14").Value = -1066
$ws.Range("H27").Value = 45670.832
$ws.Range("J27").Value = 45670.832
$ws.Range("L27").Value = 45670.832
$ws.Range("N27").Value = -45808.832
$ws.Range("H81").Value = 6199.75
$ws.Range("I81").Value = 4599
$ws.Range("K81").Value = 9198
$ws.Range("M81").Value = -8137
$ws.Range("H84").Value = 6199.75
$ws.Range("I84").Value = 4599
$ws.Range("K84").Value = 45990
$ws.Range("M84").Value = -40686
$ws.Range("H92").Value = 29000
$ws.Range("J92").Value = 29000
$ws.Range("L92").Value = 29000
$ws.Range("N92").Value = -33992
$ws.Range("H115").Value = 47924
$ws.Range("J115").Value = 47924
$ws.Range("L115").Value = 47924
$ws.Range("N115").Value = -51058
$ws.Range("H122").Value = 2996.9
$ws.Range("I122").Value = 2996.9
$ws.Range("K122").Value = 8990.700000000001
$ws.Range("M122").Value = -6540.700000000001
$ws.Range("H132").Value = 12502805
$ws.Range("I132").Value = 16130555
$ws.Range("J132").Value = 7221.5557
$ws.Range("K132").Value = 48391665
$ws.Range("L132").Value = 21664.6671
$ws.Range("M132").Value = -48389135
$ws.Range("N132").Value = -26724.6671
$ws.Range("H136").Value = 22729690
$ws.Range("I136").Value = 23811960
$ws.Range("K136").Value = 71435880
$ws.Range("M136").Value = -71433330

